# Ian Wu Ling En updated the image file used for the PutShip sequence
# diagram: a picture shape was inserted on the (only) slide and then
# the action was undone, so the picture does not remain on the slide
# afterwards -- PowerPoint still records the insert/delete/modify in
# its co-authoring change log (addSp/delSp/modSp, "add del mod").
#
# Reproduce that lifecycle here: insert the updated diagram picture,
# then immediately remove it again so the slide ends up back at its
# original, undone state.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$pic = $s.Shapes.AddPicture("PutShipSequenceDiagram.png", $false, $true, 0, 0, $p.PageSetup.SlideWidth, $p.PageSetup.SlideHeight)

# The insert is undone immediately afterwards.
$pic.Delete()
